# PageUtility class added and implemented.
#
# Changes applied on the "ManageContactPage" sheet:
#   - D2 ("Delivery Time") changes from the text "1 hour" to the number 42.5.
#     Since "1 hour" then becomes unused, it drops out of the shared-strings
#     table automatically when the workbook is saved (the remaining strings,
#     e.g. "Rs. 100", "News", "No New News", shift down by one index - this
#     happens transparently because they are addressed by text/cell, not by
#     raw shared-string index).
#   - E2 ("Delivery Charge Limit") keeps its original text, "Rs. 100".
#   - The active selection on the sheet moves from B7 to E6.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ManageContactPage")
$ws.Activate()

# Delivery Time becomes a numeric value instead of the text "1 hour".
$ws.Range("D2").Value = 42.5

# Delivery Charge Limit keeps its text.
$ws.Range("E2").Value = "Rs. 100"

# Move the sheet's active selection to E6.
$ws.Range("E6").Select()
